$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110, pushing existing rows 110-159 down to 111-160
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new weekly data point
$ws.Range("A110").Value = 7
$ws.Range("B110").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C110").Value = 'Ñuble'
$ws.Range("D110").Value = 44460
$ws.Range("E110").Value = 16
$ws.Range("F110").Value = 100112008
$ws.Range("G110").Value = 'Coliflor'
$ws.Range("H110").Value = 'Sin especificar'
$ws.Range("I110").Value = 'Primera'
$ws.Range("J110").Value = 300
$ws.Range("K110").Value = 700
$ws.Range("L110").Value = 750
$ws.Range("M110").Value = 725
$ws.Range("N110").Value = '$/unidad'
$ws.Range("O110").Value = 'Región del Maule'
$ws.Range("P110").Value = 725
$ws.Range("Q110").Value = 1
$ws.Range("R110").Value = 'Hortaliza'
